$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.000001295275857016165
$ws.Cells.Item(2, 3).Value = 0.002571899574220771
$ws.Cells.Item(2, 4).Value = 0.1494219747398047
$ws.Cells.Item(2, 5).Value = 0.4942365360607697
$ws.Cells.Item(2, 7).Value = 0.6462317056506521

$ws.Cells.Item(3, 2).Value = 1.455362044514542
$ws.Cells.Item(3, 3).Value = 1.655778082260271
$ws.Cells.Item(3, 4).Value = 0.7527432677738641
$ws.Cells.Item(3, 5).Value = 0.4942365360607697
$ws.Cells.Item(3, 7).Value = 4.358119930609447

$ws.Cells.Item(4, 2).Value = 0.0000001026826723649066
$ws.Cells.Item(4, 3).Value = 0.00006240767534437808
$ws.Cells.Item(4, 4).Value = 0.7527432677738641
$ws.Cells.Item(4, 5).Value = 0.4942365360607697
$ws.Cells.Item(4, 7).Value = 1.247042314192651

$ws.Cells.Item(5, 2).Value = 0.2917716402565462
$ws.Cells.Item(5, 3).Value = 10.34677158129881
$ws.Cells.Item(5, 4).Value = 0.1494219747398047
$ws.Cells.Item(5, 5).Value = 10.19245300693656
$ws.Cells.Item(5, 7).Value = 20.98041820323172

$ws.Cells.Item(6, 2).Value = 1.455362044514542
$ws.Cells.Item(6, 3).Value = 0.306821227259698
$ws.Cells.Item(6, 4).Value = 261.3203778131603
$ws.Cells.Item(6, 5).Value = 1133.036916526867
$ws.Cells.Item(6, 7).Value = 1396.119477611802

$ws.Cells.Item(7, 2).Value = 1.455362044514542
$ws.Cells.Item(7, 3).Value = 1.655778082260271
$ws.Cells.Item(7, 4).Value = 261.3203778131603
$ws.Cells.Item(7, 5).Value = 10.19245300693656
$ws.Cells.Item(7, 7).Value = 274.6239709468716

$ws.Cells.Item(8, 2).Value = 0.6606524410359556
$ws.Cells.Item(8, 3).Value = 1.655778082260271
$ws.Cells.Item(8, 4).Value = 0.7527432677738641
$ws.Cells.Item(8, 5).Value = 10.19245300693656
$ws.Cells.Item(8, 7).Value = 13.26162679800665

$ws.Cells.Item(9, 2).Value = 1.455362044514542
$ws.Cells.Item(9, 3).Value = 10.34677158129881
$ws.Cells.Item(9, 4).Value = 0.7527432677738641
$ws.Cells.Item(9, 5).Value = 10.19245300693656
$ws.Cells.Item(9, 7).Value = 22.74732990052378

$ws.Cells.Item(10, 2).Value = 3.286832544864788
$ws.Cells.Item(10, 3).Value = 1.655778082260271
$ws.Cells.Item(10, 4).Value = 0.1494219747398047
$ws.Cells.Item(10, 5).Value = 10.19245300693656
$ws.Cells.Item(10, 7).Value = 15.28448560880142
